$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates extracted from the published diff: cell reference -> new value.
# Column D values are price strings (may look numeric, e.g. "9.72", or contain
# multiple dots, e.g. "68.099.89") and must remain stored as text, matching the
# original t="inlineStr" cells, so we force text format before assigning them and
# restore the default "Normal" style afterwards to avoid leaving a new format
# applied to the cell.
$updates = @(
    @{ Cell = "D2"; Value = '68.099.89' }
    @{ Cell = "E2"; Value = '  +0.45%  ' }
    @{ Cell = "D3"; Value = '3.795.68' }
    @{ Cell = "E3"; Value = '  -0.38%  ' }
    @{ Cell = "E4"; Value = '  +0.10%  ' }
    @{ Cell = "D5"; Value = '601.60' }
    @{ Cell = "E5"; Value = '  +0.58%  ' }
    @{ Cell = "D6"; Value = '165.04' }
    @{ Cell = "E6"; Value = '  -1.60%  ' }
    @{ Cell = "E7"; Value = '  +0.04%  ' }
    @{ Cell = "E8"; Value = '  -0.58%  ' }
    @{ Cell = "E9"; Value = '  -1.14%  ' }
    @{ Cell = "E10"; Value = '  +0.29%  ' }
    @{ Cell = "E11"; Value = '  +2.83%  ' }
    @{ Cell = "E12"; Value = '  -1.91%  ' }
    @{ Cell = "D13"; Value = '35.81' }
    @{ Cell = "E13"; Value = '  -0.73%  ' }
    @{ Cell = "D14"; Value = '4.432.14' }
    @{ Cell = "E14"; Value = '  -0.33%  ' }
    @{ Cell = "D15"; Value = '3.771.28' }
    @{ Cell = "E15"; Value = '  -1.55%  ' }
    @{ Cell = "D16"; Value = '68.101.65' }
    @{ Cell = "E16"; Value = '  +0.45%  ' }
    @{ Cell = "D17"; Value = '18.39' }
    @{ Cell = "E17"; Value = '  -1.26%  ' }
    @{ Cell = "E18"; Value = '  +2.40%  ' }
    @{ Cell = "D19"; Value = '7.09' }
    @{ Cell = "E19"; Value = '  -0.25%  ' }
    @{ Cell = "D20"; Value = '461.31' }
    @{ Cell = "E20"; Value = '  -0.20%  ' }
    @{ Cell = "D21"; Value = '9.72' }
    @{ Cell = "D22"; Value = '0.702' }
    @{ Cell = "E22"; Value = '  -0.03%  ' }
    @{ Cell = "E23"; Value = '  -4.25%  ' }
    @{ Cell = "D24"; Value = '83.12' }
    @{ Cell = "E24"; Value = '  -0.61%  ' }
    @{ Cell = "D25"; Value = '12.03' }
    @{ Cell = "E25"; Value = '  -0.79%  ' }
    @{ Cell = "D26"; Value = '2.12' }
    @{ Cell = "E26"; Value = '  +0.15%  ' }
    @{ Cell = "D27"; Value = '10.02' }
    @{ Cell = "E27"; Value = '  -0.10%  ' }
    @{ Cell = "E28"; Value = '  -0.64%  ' }
    @{ Cell = "D29"; Value = '3.944.61' }
    @{ Cell = "E29"; Value = '  -0.29%  ' }
    @{ Cell = "E30"; Value = '  -5.28%  ' }
    @{ Cell = "B31"; Value = 'NEARProtocol' }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' }
    @{ Cell = "D31"; Value = '7.36' }
    @{ Cell = "E31"; Value = '  +0.87%  ' }
    @{ Cell = "B32"; Value = 'ImmutableX' }
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = "D32"; Value = '2.23' }
    @{ Cell = "E32"; Value = '  -0.88%  ' }
    @{ Cell = "D33"; Value = '29.36' }
    @{ Cell = "E33"; Value = '  -1.32%  ' }
    @{ Cell = "E34"; Value = '  +0.01%  ' }
    @{ Cell = "E35"; Value = '  -0.77%  ' }
    @{ Cell = "D36"; Value = '0.0998' }
    @{ Cell = "E36"; Value = '  -0.38%  ' }
    @{ Cell = "B37"; Value = 'dogwifhat' }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' }
    @{ Cell = "D37"; Value = '3.32' }
    @{ Cell = "E37"; Value = '  -3.21%  ' }
    @{ Cell = "B38"; Value = 'Kaspa' }
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = "D38"; Value = '0.139' }
    @{ Cell = "E38"; Value = '  +0.61%  ' }
    @{ Cell = "E39"; Value = '  +0.88%  ' }
    @{ Cell = "D40"; Value = '0.989' }
    @{ Cell = "E40"; Value = '  -1.34%  ' }
    @{ Cell = "D41"; Value = '0.999' }
    @{ Cell = "E41"; Value = '  -0.02%  ' }
    @{ Cell = "E43"; Value = '  -0.07%  ' }
    @{ Cell = "D44"; Value = '47.51' }
    @{ Cell = "E44"; Value = '  -1.37%  ' }
    @{ Cell = "D45"; Value = '43.31' }
    @{ Cell = "E45"; Value = '  -1.03%  ' }
    @{ Cell = "D46"; Value = '151.59' }
    @{ Cell = "E46"; Value = '  +1.60%  ' }
    @{ Cell = "E47"; Value = '  +0.32%  ' }
    @{ Cell = "B48"; Value = 'Stacks' }
    @{ Cell = "C48"; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = "D48"; Value = '1.87' }
    @{ Cell = "E48"; Value = '  +1.47%  ' }
    @{ Cell = "B49"; Value = 'ONDO' }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' }
    @{ Cell = "D49"; Value = '1.35' }
    @{ Cell = "E49"; Value = '  +2.74%  ' }
    @{ Cell = "B50"; Value = 'Bittensor' }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = "D50"; Value = '393.58' }
    @{ Cell = "E50"; Value = '  -1.41%  ' }
    @{ Cell = "D51"; Value = '26.64' }
    @{ Cell = "E51"; Value = '  -0.42%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell -match "^D") {
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}

